$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 92, shifting existing rows 92-148 down to 93-149
$ws.Rows("92:92").Insert()

# Populate the new row 92 with the new record's data (copy of row pattern, new values per diff)
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = 45086
$ws.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = 100112031
$ws.Cells.Item(92, 7).Value = "Poroto verde"
$ws.Cells.Item(92, 8).Value = "Magnum"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 45
$ws.Cells.Item(92, 11).Value = 25000
$ws.Cells.Item(92, 12).Value = 25000
$ws.Cells.Item(92, 13).Value = 25000
$ws.Cells.Item(92, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(92, 15).Value = "Perú"
$ws.Cells.Item(92, 16).Value = 1000
$ws.Cells.Item(92, 17).Value = 25
$ws.Cells.Item(92, 18).Value = "Hortaliza"
